$d = $word.ActiveDocument

# The diary table is the first (and only) table in the document. The last
# existing row records the 23/11/2021 "input recognition system" entry; we
# need to append a brand-new row underneath it describing the 27/11/2021
# "endpaint" work.
$t = $d.Tables.Item(1)

# Adding a row duplicates the formatting (cell widths, borders, shading,
# paragraph properties) of the table's current last row, which already
# matches the look of the row we need to add.
$newRow = $t.Rows.Add()
$i = $newRow.Index

$t.Cell($i, 1).Range.Text = "27/11/2021"
$t.Cell($i, 2).Range.Text = "1 Hour 30 Minutes"
$t.Cell($i, 3).Range.Text = "Stage 3"
$t.Cell($i, 4).Range.Text = "Added the endpaint function, accessed via the -kHDC command. This system call is only utilised if the user requests it, however they are given a warning if they attempt to exit the painter program space without closing the hdc. This is done to demonstrate the ability of the system to hold onto hdcs and prevent interactions between them – which is best shown through the -gHDC command which allows a user to select any hdc index desired. This is designed as solely an administrative feature that would not be included in a released version of this project, but is a necessary feature to show the working of stage 3."
